$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force a cell to hold an exact literal string (no Excel auto number/date coercion),
# while keeping its style index identical to the original (no NumberFormat bleed-through).
function Set-TextCell {
    param($cell, [string]$text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell $ws.Range("D2") '42.766.49'
Set-TextCell $ws.Range("E2") '  +0.63%  '

# Row 3
Set-TextCell $ws.Range("D3") '2.290.29'
Set-TextCell $ws.Range("E3") '  -0.17%  '

# Row 4
Set-TextCell $ws.Range("D4") '1.00'
Set-TextCell $ws.Range("E4") '  +0.10%  '

# Row 5
Set-TextCell $ws.Range("D5") '301.08'
Set-TextCell $ws.Range("E5") '  +0.13%  '

# Row 6
Set-TextCell $ws.Range("D6") '99.11'
Set-TextCell $ws.Range("E6") '  +2.26%  '

# Row 7
Set-TextCell $ws.Range("E7") '  +0.19%  '

# Row 8
Set-TextCell $ws.Range("E8") '  +0.04%  '

# Row 9
Set-TextCell $ws.Range("D9") '0.512'
Set-TextCell $ws.Range("E9") '  +3.81%  '

# Row 10
Set-TextCell $ws.Range("D10") '35.94'
Set-TextCell $ws.Range("E10") '  +7.14%  '

# Row 11
Set-TextCell $ws.Range("E11") '  -0.31%  '

# Row 12
Set-TextCell $ws.Range("E12") '  +1.92%  '

# Row 13
Set-TextCell $ws.Range("D13") '17.86'
Set-TextCell $ws.Range("E13") '  +10.98%  '

# Row 14
Set-TextCell $ws.Range("D14") '6.82'
Set-TextCell $ws.Range("E14") '  +1.34%  '

# Row 15
Set-TextCell $ws.Range("D15") '2.645.33'
Set-TextCell $ws.Range("E15") '  +0.30%  '

# Row 16
Set-TextCell $ws.Range("D16") '2.296.76'
Set-TextCell $ws.Range("E16") '  +0.67%  '

# Row 17
Set-TextCell $ws.Range("E17") '  +0.53%  '

# Row 18
Set-TextCell $ws.Range("D18") '42.672.64'
Set-TextCell $ws.Range("E18") '  +0.69%  '

# Row 19
Set-TextCell $ws.Range("D19") '12.37'
Set-TextCell $ws.Range("E19") '  +5.34%  '

# Row 20
Set-TextCell $ws.Range("E20") '  +2.66%  '

# Row 21
Set-TextCell $ws.Range("D21") '0.0₃0900'
Set-TextCell $ws.Range("E21") '  +0.28%  '

# Row 22
Set-TextCell $ws.Range("D22") '67.83'
Set-TextCell $ws.Range("E22") '  +1.72%  '

# Row 23
Set-TextCell $ws.Range("D23") '235.65'
Set-TextCell $ws.Range("E23") '  -0.16%  '

# Row 24
Set-TextCell $ws.Range("D24") '2.22'
Set-TextCell $ws.Range("E24") '  +12.73%  '

# Row 25
Set-TextCell $ws.Range("E25") '  -0.06%  '

# Row 26
Set-TextCell $ws.Range("E26") '  -0.83%  '

# Row 27
Set-TextCell $ws.Range("D27") '24.59'
Set-TextCell $ws.Range("E27") '  +2.56%  '

# Row 28
Set-TextCell $ws.Range("D28") '168.41'
Set-TextCell $ws.Range("E28") '  +0.45%  '

# Row 29
Set-TextCell $ws.Range("D29") '34.48'
Set-TextCell $ws.Range("E29") '  +1.30%  '

# Row 30
Set-TextCell $ws.Range("E30") '  -5.04%  '

# Row 31
Set-TextCell $ws.Range("E31") '  -0.37%  '

# Row 32
Set-TextCell $ws.Range("E32") '  +0.13%  '

# Row 33
Set-TextCell $ws.Range("D33") '4.97'
Set-TextCell $ws.Range("E33") '  +0.42%  '

# Row 34
Set-TextCell $ws.Range("D34") '17.51'
Set-TextCell $ws.Range("E34") '  +3.19%  '

# Row 35
Set-TextCell $ws.Range("D35") '4.59'
Set-TextCell $ws.Range("E35") '  -2.74%  '

# Row 36
Set-TextCell $ws.Range("E36") '  +2.89%  '

# Row 37
Set-TextCell $ws.Range("E37") '  -1.16%  '

# Row 38
Set-TextCell $ws.Range("B38") 'Kaspa'
Set-TextCell $ws.Range("C38") 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextCell $ws.Range("D38") '0.101'
Set-TextCell $ws.Range("E38") '  +1.92%  '

# Row 39
Set-TextCell $ws.Range("B39") 'LidoDAOToken'
Set-TextCell $ws.Range("C39") 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextCell $ws.Range("D39") '2.82'
Set-TextCell $ws.Range("E39") '  +0.26%  '

# Row 40
Set-TextCell $ws.Range("E40") '  +1.27%  '

# Row 41
Set-TextCell $ws.Range("E41") '  -0.14%  '

# Row 42
Set-TextCell $ws.Range("D42") '1.992.59'
Set-TextCell $ws.Range("E42") '  +1.40%  '

# Row 43
Set-TextCell $ws.Range("E43") '  +2.51%  '

# Row 44
Set-TextCell $ws.Range("E44") '  -2.39%  '

# Row 45
Set-TextCell $ws.Range("D45") '10.11'
Set-TextCell $ws.Range("E45") '  +4.63%  '

# Row 46
Set-TextCell $ws.Range("B46") 'EnergySwap'
Set-TextCell $ws.Range("C46") 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell $ws.Range("D46") '17.55'
Set-TextCell $ws.Range("E46") '  -0.96%  '

# Row 47
Set-TextCell $ws.Range("B47") 'NEARProtocol'
Set-TextCell $ws.Range("C47") 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell $ws.Range("D47") '2.88'
Set-TextCell $ws.Range("E47") '  +2.09%  '

# Row 48
Set-TextCell $ws.Range("D48") '55.67'
Set-TextCell $ws.Range("E48") '  +5.70%  '

# Row 49
Set-TextCell $ws.Range("D49") '2.512.75'
Set-TextCell $ws.Range("E49") '  +0.02%  '

# Row 50
Set-TextCell $ws.Range("E50") '  +2.31%  '

# Row 51
Set-TextCell $ws.Range("D51") '4.50'
Set-TextCell $ws.Range("E51") '  -0.93%  '
